$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")
$ws.Activate()

# Row 17: new work-log entry for 8.10.2018, 17:15-19:15 (2h), sprint 2.
$ws.Cells.Item(17, 1).Value = 43381
$ws.Cells.Item(17, 2).Value = 0.71875
$ws.Cells.Item(17, 3).Value = 0.80208333333333337
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = "1 h JCE-tiedoston asennusta ja etsimistä, 10 min esityslistan tekoa ja työaikakirjanpidon raportointia -> löytyvät moodlesta sekä kokouksen osanottajien sähköposteista. 50 min Windows 10 VM:n asennusta ja konffailua. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%208.10.2018.txt"

# Row grew taller to fit the wrapped text of the new note.
$ws.Rows.Item(17).RowHeight = 105

# Update the on-screen selection/scroll position to reflect where the user
# ended up after typing the new row.
$ws.Range("E17").Select()

$wb.Save()
